$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.060.56"
$ws.Range("E2").Value = "  -2.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.858.52"
$ws.Range("E3").Value = "  -3.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.12"
$ws.Range("E5").Value = "  -3.43%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4653"
$ws.Range("E7").Value = "  -2.68%  "
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06524"
$ws.Range("E9").Value = "  -3.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.43"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07808"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.22"
$ws.Range("E12").Value = "  -7.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.855.68"
$ws.Range("E13").Value = "  -3.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.100"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6623"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.89"
$ws.Range("E16").Value = "  -4.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.086.04"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.466"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.52"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.110.01"
$ws.Range("E21").Value = "  -3.28%  "
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("B23").Value = "ShibaInu"
$ws.Range("C23").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.000007203"
$ws.Range("E23").Value = "  -5.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.104"
$ws.Range("E24").Value = "  -4.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.275"
$ws.Range("E25").Value = "  -2.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.36"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.77"
$ws.Range("E27").Value = "  -4.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.897"
$ws.Range("E28").Value = "  -9.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.330"
$ws.Range("E29").Value = "  -4.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09516"
$ws.Range("E30").Value = "  -5.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.405"
$ws.Range("E31").Value = "  -4.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.464"
$ws.Range("E32").Value = "  -4.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.069"
$ws.Range("E33").Value = "  -6.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04614"
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6975"
$ws.Range("E35").Value = "  -5.34%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.092"
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.696"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01841"
$ws.Range("E38").Value = "  -5.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.255"
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.506"
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.10"
$ws.Range("E41").Value = "  -5.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8522"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.904"
$ws.Range("E44").Value = "  -5.98%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.78"
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4130"
$ws.Range("E46").Value = "  -5.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "991.46"
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.300"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.151"
$ws.Range("E49").Value = "  -4.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.87"
$ws.Range("E50").Value = "  -3.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1131"
$ws.Range("E51").Value = "  -6.48%  "
